# Auto-generated script applying scheduled-runner price/profit updates
# to the Golem_Profits workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 134999.5
$ws.Range("I51").Value = 9999.5
$ws.Range("K51").Value = 9999.5
$ws.Range("M51").Value = -9515.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1479.8
$ws.Range("I100").Value = 1479.8
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1479.8
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -938.8
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4998.6665
$ws.Range("I102").Value = 5248
$ws.Range("K102").Value = 5248
$ws.Range("M102").Value = -3626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3042.2
$ws.Range("I110").Value = 3052.75
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 3052.75
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = -1007.75
$ws.Range("N110").Value = -7090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105:N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16:N16").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 913161.75
$ws.Range("I99").Value = 836295
$ws.Range("J99").Value = 1005401.8
$ws.Range("K99").Value = 836295
$ws.Range("L99").Value = 1005401.8
$ws.Range("M99").Value = -834797
$ws.Range("N99").Value = -1008397.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113:N113").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4000
$ws.Range("J122").Value = 4000
$ws.Range("L122").Value = 12000
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 913161.75
$ws.Range("I126").Value = 836295
$ws.Range("J126").Value = 1005401.8
$ws.Range("K126").Value = 2508885
$ws.Range("L126").Value = 3016205.4
$ws.Range("M126").Value = -2506415
$ws.Range("N126").Value = -3021145.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3731.5
$ws.Range("I132").Value = 2012
$ws.Range("K132").Value = 6036
$ws.Range("M132").Value = -3506

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 252.875
$ws.Range("J17").Value = 304.16666
$ws.Range("L17").Value = 912.4999799999999
$ws.Range("N17").Value = -1250.49998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 159499.5
$ws.Range("J5").Value = 159499.5
$ws.Range("L5").Value = 159499.5
$ws.Range("N5").Value = -159725.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1067.125
$ws.Range("I16").Value = 1067.125
$ws.Range("K16").Value = 1067.125
$ws.Range("M16").Value = -897.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 738.25
$ws.Range("I22").Value = 691.1429000000001
$ws.Range("J22").Value = 774.8889
$ws.Range("K22").Value = 691.1429000000001
$ws.Range("L22").Value = 774.8889
$ws.Range("M22").Value = -396.1429000000001
$ws.Range("N22").Value = -1364.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 738.25
$ws.Range("I27").Value = 691.1429000000001
$ws.Range("J27").Value = 774.8889
$ws.Range("K27").Value = 691.1429000000001
$ws.Range("L27").Value = 774.8889
$ws.Range("M27").Value = -584.1429000000001
$ws.Range("N27").Value = -988.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11724.1
$ws.Range("I61").Value = 1779.5
$ws.Range("K61").Value = 1779.5
$ws.Range("M61").Value = -1577.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3951
$ws.Range("I68").Value = 5902
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 5902
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -5153
$ws.Range("N68").Value = -3498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3951
$ws.Range("I71").Value = 5902
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 29510
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -25766
$ws.Range("N71").Value = -17488

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 495.18182
$ws.Range("I82").Value = 411
$ws.Range("K82").Value = 411
$ws.Range("M82").Value = -50

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 495.18182
$ws.Range("I85").Value = 411
$ws.Range("K85").Value = 411
$ws.Range("M85").Value = 837

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 11724.1
$ws.Range("I113").Value = 1779.5
$ws.Range("K113").Value = 1779.5
$ws.Range("M113").Value = 390.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 950
$ws.Range("I81").Value = 950
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1900
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -839
$ws.Range("N81").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 950
$ws.Range("I84").Value = 950
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9500
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4196
$ws.Range("N84").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1600
$ws.Range("I107").Value = 1200
$ws.Range("K107").Value = 3600
$ws.Range("M107").Value = -1680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2993
$ws.Range("I122").Value = 2632
$ws.Range("K122").Value = 7896
$ws.Range("M122").Value = -5446
